$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B5").Value = "Alleine an Dartsturnier"
$ws.Range("C5").Value = "Ich habe mich gestern für ein Darts Turnier am Sonntag angemeldet, obwohl mein Bruder nicht kann. Somit gehe ich alleine und lerne neue Menschen kennen."
$ws.Range("D5").Value = "Mut"
$ws.Range("E5").Value = "https://as1.ftcdn.net/v2/jpg/00/36/21/18/1000_F_36211854_RIs7qAGng1K4c5JvgxN3HxxMHbFtLMlc.jpg"

$ws.Range("E7").Select()
